$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted above the current row 6,
# pushing all existing data rows (old 6..81) down by one (new 7..82).
$ws.Rows(6).Insert()

# Populate the newly inserted row 6 with this week's record.
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44750
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100108
$ws.Range("H6").Value = "Tropicales y subtropicales"
$ws.Range("I6").Value = 100108007
$ws.Range("J6").Value = "Coco"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 28000
$ws.Range("O6").Value = 28000
$ws.Range("P6").Value = 28000
$ws.Range("Q6").Value = "$/malla 20 unidades"
$ws.Range("R6").Value = "Perú"
$ws.Range("S6").Value = 1400
$ws.Range("T6").Value = 20
